# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2210"
#   "<name>_new" -> "<name>_FV2304"
# Then wrap the sheet's used range in a native Excel Table ("Table1"),
# and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row (row 1) cells -----------------------------
$lastCol = 21   # A..U
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cur = $cell.Value2
    if ($cur -ne $null) {
        if ($cur -like "*_old") {
            $cell.Value2 = ($cur -replace "_old$", "_FV2210")
        } elseif ($cur -like "*_new") {
            $cell.Value2 = ($cur -replace "_new$", "_FV2304")
        }
    }
}

# --- 2) Turn A1:U64 into a real Table (ListObject) ----------------------
# Adding a ListObject directly on top of the already-formatted header row
# makes Excel capture a "header row" dxf (bold/fill/border) override, and
# also forces a named table style - neither of which is present in the
# target workbook. So we build the table on a blank scratch row first
# (where the header cells carry no special formatting) and then move it
# back onto the real header row with Resize - this keeps the original
# cell styles/style table untouched and produces a plain table definition.
$scratchRow = 1000
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item($scratchRow, $c).Value2 = $ws.Cells.Item(1, $c).Value2
}

$scratchRange = $ws.Range($ws.Cells.Item($scratchRow, 1), $ws.Cells.Item($scratchRow, $lastCol))
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $scratchRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# clear the scratch header text, then move the table over the real data
$scratchRange.ClearContents()
$fullRange = $ws.Range("A1:U64")
$tbl.Resize($fullRange)

# --- 3) Freeze the header row --------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
